# Présentation du projet 4 — add a new slide before the final
# "Fin de la présentation" slide, containing the list of changes made
# according to the accessibility checklist.

$p = $ppt.ActivePresentation

# The original last slide (29) is "Fin de la présentation" — keep it
# untouched, but insert a brand-new slide right before it using the
# same "Titre et contenu" layout (same layout as slide 28/29).
$layout = $p.Slides.Item(28).CustomLayout
$newSlide = $p.Slides.AddSlide(29, $layout)

# --- Title placeholder -------------------------------------------------
$title = $newSlide.Shapes.Item(1)
$title.Name = "Titre 1"
$titleTr = $title.TextFrame.TextRange
$titleTr.Text = "Changements effectués en fonction de la checklist"
$titleTr.LanguageID = "fr-FR"

# --- Content placeholder -------------------------------------------------
$content = $newSlide.Shapes.Item(2)
$content.Name = "Espace réservé du contenu 2"
$body = $content.TextFrame.TextRange

$body.Text = "Structuration du site à l’aide de balises h"
$body.LanguageID = "fr-FR"

$p2 = $body.InsertAfter("`rAgrandissement de la zone tactile de certains boutons et liens")
$p2.LanguageID = "fr-FR"

$p3 = $body.InsertAfter("`rPrésence d’un label dans chaque champ du formulaire")
$p3.LanguageID = "fr-FR"

$p4 = $body.InsertAfter("`rTexte alternatif pour les liens")
$p4.LanguageID = "fr-FR"

$p5a = $body.InsertAfter("`rAmélioration de la visibilité en supprimant les ")
$p5a.LanguageID = "fr-FR"

$p5b = $body.InsertAfter("z-index,display")
$p5b.LanguageID = "fr-FR"

$p5c = $body.InsertAfter(" inutiles")
$p5c.LanguageID = "fr-FR"

$p6 = $body.InsertAfter("`rAugmentation de la taille de la police si inférieur à 12px, ")
$p6.LanguageID = "fr-FR"

$p7a = $body.InsertAfter("`rAugmentation des contrastes entre l’arrière plan et les éléments de ")
$p7a.LanguageID = "fr-FR"

$p7b = $body.InsertAfter("premier plan")
$p7b.LanguageID = "fr-FR"

$p8 = $body.InsertAfter("`r`r")
$p8.LanguageID = "fr-FR"

Write-Output "Slides after edit:"
Write-Output $p.Slides.Count
